$wb = $excel.ActiveWorkbook

# --- Re-create the "总计" (total) sheet so the new "2022-Q1" sheet can slot in
#     between "2021-Q3" and "总计" with the expected sheetId sequence ---
$total = $wb.Worksheets.Item("总计")
$total.Delete() | Out-Null

$q1After = $wb.Worksheets.Item("2021-Q3")
$wsQ1 = $wb.Worksheets.Add($null, $q1After)
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

# ---------------------------------------------------------------------------
# "2022-Q1" sheet data
# ---------------------------------------------------------------------------
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

$wsQ1.Range("A2").Value = 0
$wsQ1.Range("B2").Value = "'009387"
$wsQ1.Range("C2").Value = "嘉实稳福混合A"
$wsQ1.Range("D2").Value = "'0.08"
$wsQ1.Range("E2").Value = "'34.71"
$wsQ1.Range("F2").Value = "'1.59"
$wsQ1.Range("G2").Value = "'0.0013"
$wsQ1.Range("H2").Value = 9

$wsQ1.Range("A3").Value = 1
$wsQ1.Range("B3").Value = "'009388"
$wsQ1.Range("C3").Value = "嘉实稳福混合C"
$wsQ1.Range("D3").Value = "'0.01"
$wsQ1.Range("E3").Value = "'34.71"
$wsQ1.Range("F3").Value = "'1.59"
$wsQ1.Range("G3").Value = "'0.0002"
$wsQ1.Range("H3").Value = 9

$hdrQ1 = $wsQ1.Range("B1:H1")
$hdrQ1.Font.Bold = $true
$hdrQ1.HorizontalAlignment = -4108
$hdrQ1.VerticalAlignment = -4160
$hdrQ1.Borders.LineStyle = 1

$idxQ1 = $wsQ1.Range("A2:A3")
$idxQ1.Font.Bold = $true
$idxQ1.HorizontalAlignment = -4108
$idxQ1.VerticalAlignment = -4160
$idxQ1.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# "总计" sheet data
# ---------------------------------------------------------------------------
$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q3"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 0.24

$hdrTotal = $wsTotal.Range("B1:D1")
$hdrTotal.Font.Bold = $true
$hdrTotal.HorizontalAlignment = -4108
$hdrTotal.VerticalAlignment = -4160
$hdrTotal.Borders.LineStyle = 1

$idxTotal = $wsTotal.Range("A2:A3")
$idxTotal.Font.Bold = $true
$idxTotal.HorizontalAlignment = -4108
$idxTotal.VerticalAlignment = -4160
$idxTotal.Borders.LineStyle = 1

$q1After.Select() | Out-Null
